$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @{
    2 = "6/18/2018 00:07:00"
    3 = "6/18/2018 00:08:00"
    4 = "6/18/2018 00:36:00"
    5 = "6/18/2018 00:36:00"
    6 = "6/18/2018 00:36:00"
    7 = "6/18/2018 00:37:00"
    8 = "6/18/2018 00:37:00"
    9 = "6/18/2018 01:01:00"
    10 = "6/18/2018 01:02:00"
    11 = "6/18/2018 01:02:00"
    12 = "6/18/2018 01:15:00"
    13 = "6/18/2018 01:15:00"
    14 = "6/18/2018 01:15:00"
    15 = "6/18/2018 01:16:00"
    16 = "6/18/2018 01:16:00"
    17 = "6/18/2018 01:16:00"
    18 = "6/18/2018 01:17:00"
    19 = "6/18/2018 01:17:00"
    20 = "6/18/2018 01:17:00"
    21 = "6/18/2018 01:17:00"
    22 = "6/18/2018 01:17:00"
    23 = "6/18/2018 01:18:00"
    24 = "6/18/2018 01:18:00"
    25 = "6/18/2018 01:23:00"
    26 = "6/18/2018 01:23:00"
    27 = "6/18/2018 01:24:00"
    28 = "6/18/2018 01:24:00"
    29 = "6/18/2018 01:24:00"
    30 = "6/18/2018 01:24:00"
    31 = "6/18/2018 01:30:00"
    32 = "6/18/2018 01:30:00"
    33 = "6/18/2018 01:30:00"
    34 = "6/18/2018 01:30:00"
    35 = "6/18/2018 01:30:00"
    36 = "6/18/2018 01:30:00"
    37 = "6/18/2018 01:31:00"
    38 = "6/18/2018 01:31:00"
    39 = "6/18/2018 01:36:00"
    40 = "6/18/2018 01:36:00"
    41 = "6/18/2018 01:36:00"
    42 = "6/18/2018 01:36:00"
    43 = "6/18/2018 01:36:00"
    44 = "6/18/2018 01:37:00"
    45 = "6/18/2018 01:37:00"
    46 = "6/18/2018 01:37:00"
    47 = "6/18/2018 01:37:00"
    48 = "6/18/2018 01:37:00"
    49 = "6/18/2018 01:37:00"
    50 = "6/18/2018 01:37:00"
    51 = "6/18/2018 01:37:00"
    52 = "6/18/2018 01:37:00"
    53 = "6/18/2018 01:38:00"
    54 = "6/18/2018 01:39:00"
    55 = "6/18/2018 01:39:00"
    56 = "6/18/2018 01:39:00"
    57 = "6/18/2018 01:39:00"
    58 = "6/18/2018 01:40:00"
    59 = "6/18/2018 01:40:00"
    60 = "6/18/2018 01:40:00"
    61 = "6/18/2018 01:41:00"
    62 = "6/18/2018 01:41:00"
    63 = "6/18/2018 01:42:00"
    64 = "6/18/2018 01:42:00"
    65 = "6/18/2018 09:18:00"
    66 = "6/18/2018 09:18:00"
    67 = "6/18/2018 09:18:00"
    68 = "6/18/2018 09:19:00"
    69 = "6/18/2018 09:19:00"
    70 = "6/18/2018 09:21:00"
    71 = "6/18/2018 09:22:00"
    72 = "6/18/2018 09:22:00"
    73 = "6/18/2018 09:23:00"
    74 = "6/18/2018 09:23:00"
    75 = "6/18/2018 09:24:00"
    76 = "6/18/2018 09:24:00"
    77 = "6/18/2018 09:24:00"
    78 = "6/18/2018 09:24:00"
    79 = "6/18/2018 09:24:00"
    80 = "6/18/2018 09:24:00"
    81 = "6/18/2018 09:25:00"
    82 = "6/18/2018 09:25:00"
    83 = "6/18/2018 09:25:00"
    84 = "6/18/2018 09:27:00"
    85 = "6/18/2018 09:28:00"
    86 = "6/18/2018 09:34:00"
    87 = "6/18/2018 09:34:00"
    88 = "6/18/2018 09:35:00"
    89 = "6/18/2018 09:35:00"
    90 = "6/18/2018 09:35:00"
    91 = "6/18/2018 09:35:00"
    92 = "6/18/2018 09:35:00"
    93 = "6/18/2018 09:43:00"
    94 = "6/18/2018 09:44:00"
    95 = "6/18/2018 09:44:00"
    96 = "6/18/2018 09:45:00"
    97 = "6/18/2018 09:45:00"
    98 = "6/18/2018 09:53:00"
    99 = "6/18/2018 09:55:00"
    100 = "6/18/2018 09:55:00"
    101 = "6/18/2018 09:55:00"
    102 = "6/18/2018 09:55:00"
    103 = "6/18/2018 09:55:00"
    104 = "6/18/2018 09:55:00"
    105 = "6/18/2018 09:55:00"
    106 = "6/18/2018 09:56:00"
    107 = "6/18/2018 09:56:00"
    108 = "6/18/2018 09:56:00"
    109 = "6/18/2018 09:56:00"
    110 = "6/18/2018 09:56:00"
    111 = "6/18/2018 09:56:00"
    112 = "6/18/2018 09:56:00"
    113 = "6/18/2018 09:56:00"
    114 = "6/18/2018 09:57:00"
    115 = "6/18/2018 11:23:00"
    116 = "6/18/2018 11:24:00"
    117 = "6/18/2018 11:24:00"
    118 = "6/18/2018 11:35:00"
    119 = "6/18/2018 11:36:00"
    120 = "10/29/2018 11:59:00"
    121 = "10/29/2018 12:26:00"
    122 = "10/30/2018 15:55:00"
    123 = "10/30/2018 15:55:00"
    124 = "10/30/2018 15:55:00"
    125 = "10/30/2018 15:55:00"
    126 = "10/30/2018 15:55:00"
    127 = "10/30/2018 15:56:00"
    128 = "10/30/2018 15:56:00"
    129 = "10/30/2018 15:56:00"
    130 = "10/30/2018 15:56:00"
    131 = "10/30/2018 15:56:00"
    132 = "10/30/2018 15:56:00"
    133 = "10/30/2018 15:56:00"
    134 = "10/30/2018 15:56:00"
    135 = "10/30/2018 15:56:00"
    136 = "10/30/2018 15:57:00"
    137 = "10/30/2018 15:57:00"
    138 = "10/30/2018 15:57:00"
    139 = "10/30/2018 15:57:00"
    140 = "10/30/2018 15:57:00"
    141 = "10/30/2018 15:57:00"
    142 = "10/30/2018 15:57:00"
    143 = "10/30/2018 15:57:00"
    144 = "10/30/2018 15:57:00"
    145 = "11/8/2018 14:04:00"
    146 = "11/8/2018 14:05:00"
    147 = "11/12/2018 12:31:00"
    148 = "11/12/2018 12:32:00"
    149 = "11/12/2018 12:32:00"
    150 = "1/29/2019 16:38:51"
    151 = "8/22/2019 14:16:18"
    152 = "8/22/2019 14:16:22"
    153 = "8/22/2019 14:16:28"
}

foreach ($row in $dates.Keys) {
    $ws.Range("M$row").Value = $dates[$row]
}
